$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data for the "getFirstIdRoot" screen (niraqa data refresh),
# mirroring the pattern of the other *Root rows: one row per role
# (REGISTRATION_OFFICER, REGISTRATION_ADMIN, REGISTRATION_SUPERVISOR).
$ws.Range("A23").Value = "eng"
$ws.Range("B23").Value = "getFirstIdRoot"
$ws.Range("C23").Value = "REGISTRATION_OFFICER"

$ws.Range("A24").Value = "eng"
$ws.Range("B24").Value = "getFirstIdRoot"
$ws.Range("C24").Value = "REGISTRATION_ADMIN"

$ws.Range("A25").Value = "eng"
$ws.Range("B25").Value = "getFirstIdRoot"
$ws.Range("C25").Value = "REGISTRATION_SUPERVISOR"

# Columns D:E store the literal text "TRUE" (not a boolean) for every row, as
# seen elsewhere on the sheet. Typing "TRUE" directly gets auto-coerced to a
# real boolean, so copy existing "TRUE" text cells down instead (matching
# range sizes so the paste lines up 1:1) - this keeps the value/type/style
# identical to the rest of the column.
$ws.Range("D2:D4").Copy() | Out-Null
$ws.Range("D23:D25").PasteSpecial(-4104) | Out-Null
$ws.Range("E2:E4").Copy() | Out-Null
$ws.Range("E23:E25").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C25").Select()
